$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("LYFT")

$ws.Range("B6").Value = 698000000.0
$ws.Range("C6").Value = 653000000.0
$ws.Range("D6").Value = 698819000.0
$ws.Range("E6").Value = 768664000.0
$ws.Range("F6").Value = 322902000.0
